# Refresh the cryptocurrency price/volume table (rows 2-51) to the latest
# scraped values, matching the GitHub Actions "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.060.88'
$ws.Range("E2").Value = '  -0.89%  '
# Row 3
$ws.Range("D3").Value = '2.221.46'
$ws.Range("E3").Value = '  -1.01%  '
# Row 4
$ws.Range("E4").Value = '  +0.00%  '
# Row 5
$ws.Range("D5").Value = '''243.01'
$ws.Range("E5").Value = '  -0.81%  '
# Row 6
$ws.Range("D6").Value = '''0.626'
$ws.Range("E6").Value = '  -0.53%  '
# Row 7
$ws.Range("D7").Value = '''74.05'
$ws.Range("E7").Value = '  -2.23%  '
# Row 8
$ws.Range("E8").Value = '  +0.10%  '
# Row 9
$ws.Range("D9").Value = '''0.603'
$ws.Range("E9").Value = '  -3.13%  '
# Row 10
$ws.Range("D10").Value = '''42.69'
$ws.Range("E10").Value = '  -3.78%  '
# Row 11
$ws.Range("D11").Value = '''0.0961'
$ws.Range("E11").Value = '  +1.16%  '
# Row 12
$ws.Range("E12").Value = '  -3.75%  '
# Row 13
$ws.Range("D13").Value = '''0.103'
$ws.Range("E13").Value = '  +0.70%  '
# Row 14
$ws.Range("D14").Value = '2.548.96'
$ws.Range("E14").Value = '  -1.17%  '
# Row 15
$ws.Range("D15").Value = '''14.28'
$ws.Range("E15").Value = '  -2.23%  '
# Row 16
$ws.Range("D16").Value = '''0.836'
$ws.Range("E16").Value = '  -2.89%  '
# Row 17
$ws.Range("D17").Value = '2.240.13'
$ws.Range("E17").Value = '  -0.36%  '
# Row 18
$ws.Range("D18").Value = '41.887.81'
$ws.Range("E18").Value = '  -0.98%  '
# Row 19
$ws.Range("E19").Value = '  +4.09%  '
# Row 20
$ws.Range("E20").Value = '  +0.11%  '
# Row 21
$ws.Range("D21").Value = '''72.80'
$ws.Range("E21").Value = '  +0.94%  '
# Row 22
$ws.Range("D22").Value = '''11.00'
$ws.Range("E22").Value = '  -3.32%  '
# Row 23
$ws.Range("D23").Value = '''230.10'
$ws.Range("E23").Value = '  -0.77%  '
# Row 24
$ws.Range("E24").Value = '  -6.59%  '
# Row 25
$ws.Range("E25").Value = '  +0.04%  '
# Row 26
$ws.Range("D26").Value = '''11.38'
$ws.Range("E26").Value = '  -3.90%  '
# Row 27
$ws.Range("E27").Value = '  -0.08%  '
# Row 28
$ws.Range("D28").Value = '''2.28'
$ws.Range("E28").Value = '  -1.50%  '
# Row 29
$ws.Range("E29").Value = '  -3.13%  '
# Row 30
$ws.Range("D30").Value = '''166.66'
$ws.Range("E30").Value = '  -0.30%  '
# Row 31
$ws.Range("D31").Value = '''20.55'
$ws.Range("E31").Value = '  -0.91%  '
# Row 32
$ws.Range("E32").Value = '  -2.96%  '
# Row 33
$ws.Range("E33").Value = '  -2.25%  '
# Row 34
$ws.Range("D34").Value = '''30.10'
$ws.Range("E34").Value = '  -2.70%  '
# Row 35
$ws.Range("E35").Value = '  -0.84%  '
# Row 36
$ws.Range("E36").Value = '  -8.40%  '
# Row 37
$ws.Range("D37").Value = '''4.30'
$ws.Range("E37").Value = '  -8.32%  '
# Row 38
$ws.Range("E38").Value = '  -3.98%  '
# Row 39
$ws.Range("D39").Value = '''13.19'
$ws.Range("E39").Value = '  -4.86%  '
# Row 40
$ws.Range("D40").Value = '''2.12'
$ws.Range("E40").Value = '  -2.68%  '
# Row 41
$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").Value = '''5.68'
$ws.Range("E41").Value = '  -1.54%  '
# Row 42
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = '''64.79'
$ws.Range("E42").Value = '  +1.59%  '
# Row 43
$ws.Range("D43").Value = '''0.199'
$ws.Range("E43").Value = '  -1.65%  '
# Row 44
$ws.Range("D44").Value = '''8.69'
$ws.Range("E44").Value = '  -1.83%  '
# Row 45
$ws.Range("D45").Value = '''104.10'
$ws.Range("E45").Value = '  -2.78%  '
# Row 46
$ws.Range("E46").Value = '  -2.02%  '
# Row 47
$ws.Range("E47").Value = '  -3.74%  '
# Row 48
$ws.Range("E48").Value = '  -2.87%  '
# Row 49
$ws.Range("E49").Value = '  -1.05%  '
# Row 50
$ws.Range("E50").Value = '  -1.32%  '
# Row 51
$ws.Range("D51").Value = '2.424.68'
$ws.Range("E51").Value = '  -1.36%  '
